# IdentitasSarana import/export fix: add Tipe/ID/Keterangan columns,
# change row-3 "Nama Sarana" from Pel to Router, and populate the new
# Tipe/ID/Keterangan sample data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- column widths for the new columns ------------------------------------
$ws.Columns("C").ColumnWidth = 14.92
$ws.Columns("E").ColumnWidth = 9.92
$ws.Columns("F").ColumnWidth = 17.42

# --- C1 header ("Tipe") - same style as the existing header cells ---------
$ws.Range("A1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("C1").Value = "Tipe"

# --- E1:F1 headers ("ID", "Keterangan") - header style + vertical center --
$ws.Range("A1").Copy()
$ws.Range("E1:F1").PasteSpecial(-4122)
$ws.Range("E1:F1").VerticalAlignment = -4108
$ws.Range("E1").Value = "ID"
$ws.Range("F1").Value = "Keterangan"

# --- C2:C3 ("Tipe" values) - number style, centered, no wrap --------------
$ws.Range("A2").Copy()
$ws.Range("C2:C3").PasteSpecial(-4122)
$ws.Range("C2:C3").WrapText = $false
$ws.Range("C2").Value = 0
$ws.Range("C3").Value = 1

# --- E2:E3 ("ID" values) - number style, centered, vertical center --------
$ws.Range("A2").Copy()
$ws.Range("E2:E3").PasteSpecial(-4122)
$ws.Range("E2:E3").WrapText = $false
$ws.Range("E2:E3").VerticalAlignment = -4108
$ws.Range("E2").Value = 0
$ws.Range("E3").Value = 1

# --- F2:F3 ("Keterangan" values) - text style, left, vertical center ------
$ws.Range("B2").Copy()
$ws.Range("F2:F3").PasteSpecial(-4122)
$ws.Range("F2:F3").WrapText = $false
$ws.Range("F2:F3").VerticalAlignment = -4108
$ws.Range("F2").Value = "Bukan perangkat IT"
$ws.Range("F3").Value = "Perangakt IT"

# --- B3 value change: "Pel" -> "Router" ------------------------------------
$ws.Range("B3").Value = "Router"

$excel.CutCopyMode = $false

# --- match the saved selection state ---------------------------------------
[void]$ws.Range("E16").Select()
